$d = $word.ActiveDocument

$replacements = @(
    ,@("Kichwa cha Video", "Video Title")
    ,@("Tatizo la Uwanja wa Ndege", "The Airport Problem")
    ,@("Mada", "Topic")
    ,@("Jiometri", "Geometry")
    ,@("Malengo", "Aim(s)")
    ,@("Pata wazo angavu la tatizo la kupunguza, tambua jinsi ya kutekeleza kwa vitendo matatizo ya kupunguza.", "Get the intuitive idea of a minimization problem, figure out how to practically implement minimization problems.")
    ,@("Urefu", "Length")
    ,@("Mahali pa Kambi", "Camp Location")
    ,@("Wawezeshaji", "Facilitators")
    ,@("N. ya wanafunzi", "N. of students")
    ,@("Tarehe", "Date")
    ,@("Rasilimali", "Resources")
    ,@("inahitajika", "needed")
    ,@("Pini (3 kila kikundi), kamba (1/kikundi), pete ya chuma (si lazima lakini inafaa kuzuia msuguano 1/kikundi), kadibodi nene au sehemu ya mbao inayoweza kutupwa (1/kikundi)", "Pins (3 each group), string (1/group), metal ring (optional but convenient to avoid friction 1/group), thick cardboard or wooden disposable surface (1/group)")
    ,@("Maandalizi", "Preparations")
    ,@("Weka alama 3 kwenye kuni", "Pin 3 points on the wood")
    ,@("Muda wa video", "Video time")
    ,@("Mwezeshaji anafanya nini", "What facilitator does")
    ,@("Wanachofanya wanafunzi", "What learners do")
    ,@("Utangulizi Mkuu wa Video ya VMC", "General VMC Video Introduction")
    ,@("Utangulizi wa Video", "Video Introduction")
    ,@("Kitendawili", "Riddle")
    ,@("Utangulizi wa jaribio la kwanza", "Introduction of the first experiment")
    ,@("SITISHA VIDEO", "VIDEO PAUSE")
    ,@("Kutafuta suluhu", "Finding a solution")
    ,@("Kusaidia mchakato, kuchochea mawazo", "Assist the process, provoke thoughts")
    ,@("Jaribu kupata mpangilio wa kamba ili upunguzaji wa urefu wa kamba ulingane na kupunguza jumla ya urefu wa barabara", "Try to find a setting of the string such that the minimization of the length of the string corresponds to minimize the sum of the lengths of the roads")
    ,@("Suluhisho ", "Solution ")
    ,@("SITISHA VIDEO", "VIDEO PAUSE")
    ,@("Jiometri", "Geometry")
    ,@("Kusaidia mchakato, kuchochea mawazo", "Assist the process, provoke thoughts")
    ,@("Jaribu kujua ni mali gani ya kijiometri ambayo hatua mpya ina uhusiano na 3 ya kuanzia.", "Try to figure out what geometrical property the new point has in relation to the starting 3.")
    ,@("Inaonyesha pembe 120°", "Showing the 120° angles")
    ,@("suluhisho", "solution")
    ,@("Pointi ,F, inayopatikana kama inavyoonyeshwa kwenye video, inaitwa Fermat Point.", "The point ,F, found as shown in the video, is called Fermat Point.")
    ,@("Kuna njia kadhaa tofauti zinazowezekana za kamba ambazo zinaweza kutumika kupata uhakika F.", "There are several different possible paths of the string that can be used to find point F.")
    ,@("Ona kwamba matumizi ya pete sio lazima kabisa, lakini inasaidia kupunguza msuguano (adui wa asili wa uzoefu huu).", "Notice that the use of the ring is not strictly necessary, but it helps to reduce friction (natural enemy of this experience).")
    ,@("Mara tu uhakika unapopatikana (kabla ya kutazama suluhu katika video) wanafunzi wanaweza kuulizwa kutafuta pembe kwa kutambua kwamba kila pembe ina mshikamano na zinaunda 360° zote kwa pamoja.", "Once the point is found (before watching the solution in the video) students can be asked to find the angles by noticing that each angle is congruent and they form 360° all together.")
    ,@("Muundo wa kijiometri ambao unaweza kutumika kupata F ni pamoja na kujenga pembetatu zilizo sawa kwenye kando ya pembetatu asilia na kuunganisha alama tofauti:", "A geometrical construction that can be used to find F consists in building equilateral triangles on the sides of the original triangle and connecting opposite points:")
    ,@("Ujenzi huu unaweza kuigwa kwenye ubao wa mbao ili kuthibitisha kwamba njia hizo mbili zitaongoza kwenye hatua sawa.", "This construction can be replicated on the wooden board to verify that the two methods will lead to the same point.")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# Update the document default language (Swahili, Kenya -> Swahili, Tanzania)
$d.Styles("Normal").Font.LanguageID = "sw-TZ"

Write-Output "done"
